# The edit described by the diff is a cyclic re-ordering of the five data
# rows 75-79 (the row immediately under each record's original position
# gets the next record's data; one record wraps back around to the top).
# Concretely, using the *original* row numbers as keys:
#
#   old row 75  ->  new row 76
#   old row 76  ->  new row 78
#   old row 77  ->  new row 79
#   old row 78  ->  new row 77
#   old row 79  ->  new row 75
#
# Every single cell-level change in the diff (Id, Antal, coordinates,
# accuracy, times, reporter names, species fields, ...) is fully explained
# by this move: each row carries its entire original content (including
# which columns are populated) to its new location. No independent value
# edits are layered on top of the move.
#
# Because the move is a permutation (not a simple shift), the rows are
# first staged into unused far-away rows so that the overlapping
# source/destination ranges of rows 75-79 do not clobber one another,
# and only afterwards copied into their final destinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcFirst = 75
$srcLast = 79
$stageFirst = 1075
$stageLast = 1079

# 1) Stage the current content of rows 75-79 into rows 1075-1079
#    (row 75 -> 1075, 76 -> 1076, ..., 79 -> 1079).
for ($i = 0; $i -le ($srcLast - $srcFirst); $i++) {
    $srcRow = $srcFirst + $i
    $stageRow = $stageFirst + $i
    $ws.Range("A" + $srcRow + ":AY" + $srcRow).Copy($ws.Range("A" + $stageRow + ":AY" + $stageRow))
}

# 2) Clear rows 75-79 completely. This Excel/COM shim does not blank out
#    destination cells whose corresponding source cell is empty when using
#    Copy(), so the old content must be removed explicitly before pasting
#    the permuted rows back in, otherwise stale values would survive.
$ws.Range("A" + $srcFirst + ":AY" + $srcLast).ClearContents()

# 3) Copy the staged rows back to their permuted destinations.
$moveMap = @{ 75 = 76; 76 = 78; 77 = 79; 78 = 77; 79 = 75 }
foreach ($oldRow in $moveMap.Keys) {
    $newRow = $moveMap[$oldRow]
    $stageRow = $stageFirst + ($oldRow - $srcFirst)
    $ws.Range("A" + $stageRow + ":AY" + $stageRow).Copy($ws.Range("A" + $newRow + ":AY" + $newRow))
}

# 4) Remove the temporary staging rows.
$ws.Range("A" + $stageFirst + ":AY" + $stageLast).ClearContents()
